$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.876671075820923
$ws.Range("B1").Value = 3.021333932876587
$ws.Range("C1").Value = 2.737329959869385
$ws.Range("D1").Value = 3.071438789367676
$ws.Range("E1").Value = 2.599180698394775
